$d = $word.ActiveDocument

# The document currently ends with six paragraphs (right after
# "Figura 3: Sinais transmitidos."): a lone "_GoBack" bookmark paragraph,
# four empty centered spacer paragraphs, and the paragraph holding the
# second figure (Saida_Filtro.bmp). We replace each of them, one at a
# time (re-resolving the Range fresh before every call — replacing a
# multi-paragraph Range in one InsertXML call does not reliably delete
# the original paragraph marks), with the six paragraphs described by
# the target revision: the new body paragraph, the image paragraph
# (now carrying a lastRenderedPageBreak), the new figure caption, and
# three new analysis paragraphs — the last of which now carries the
# relocated "_GoBack" bookmark.

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Replace-ParagraphXml($paraIndex, $innerXml) {
    $d2 = $word.ActiveDocument
    $p = $d2.Paragraphs($paraIndex)
    $r = $d2.Range($p.Range.Start, $p.Range.End)
    $r.InsertXML($xmlHeader + $innerXml + $xmlFooter)
}

$count = $d.Paragraphs.Count
$first = $count - 5   # the lone "_GoBack" bookmark paragraph

$frag0 = '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">Com relação à atividade 2, foi obtida a figura </w:t></w:r><w:r><w:t>seguinte</w:t></w:r><w:r><w:t>, que mostra a amostragem do sinal recebido.</w:t></w:r></w:p>'
Replace-ParagraphXml $first $frag0

$frag1 = '<w:p w:rsidR="00AA48C0" w:rsidRPr="001012BC" w:rsidRDefault="00951630" w:rsidP="00951630"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="5334000" cy="4000500"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1" name="Imagem 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name="Saida_Filtro.bmp"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId10"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5334000" cy="4000500"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'
Replace-ParagraphXml ($first + 1) $frag1

$frag2 = '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Figura 4: Sinal Recebido.</w:t></w:r></w:p>'
Replace-ParagraphXml ($first + 2) $frag2

$frag3 = '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>Em termos de BER, foi observado que a BER desse sinal foi pior que o sinal do laboratório 2, pois agora há interferência entre os símbolos transmitidos.</w:t></w:r></w:p>'
Replace-ParagraphXml ($first + 3) $frag3

$frag4 = '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Com relação ao atraso, foi observado que quanto </w:t></w:r><w:r><w:t>menor</w:t></w:r><w:r><w:t xml:space="preserve"> o valor de alfa, maior era o atraso, pois a </w:t></w:r><w:r><w:t>resposta</w:t></w:r><w:r><w:t xml:space="preserve"> era maior e levava mais tempo para ocorrer o tempo de amostragem.</w:t></w:r></w:p>'
Replace-ParagraphXml ($first + 4) $frag4

$frag5 = '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>Assim, n</w:t></w:r><w:r><w:t>ão usamos na prátic</w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t xml:space="preserve"> o menor valor de alfa</w:t></w:r><w:r><w:t xml:space="preserve"> pois isso geraria mais atraso de amostragem, o que também diminuiria o desempenho do sistema.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Replace-ParagraphXml ($first + 5) $frag5

Write-Host "Final paragraph count:" $d.Paragraphs.Count
